$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the date cell format (numFmtId 14, centered) from B3 to B5:B6 first
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B5:B6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Row 5: Day 2, Date 2024-06-06, Time Spent 1, Description "Am terminat Obiective + Structura Lucrarii"
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 45449
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = "Am terminat Obiective + Structura Lucrarii"

# Row 6: Day 3, Date 2024-07-06
$ws.Range("A6").Value = 3
$ws.Range("B6").Value = 45479

# Update the active selection to match the target state
$ws.Range("C8").Select()
